# Add Employee Related changes
# - TestCases sheet: record pass/fail Result for two more test rows (G3, G7)
# - Add Employee sheet: insert UserName/Password login columns before the
#   existing FirstName/MiddleName/LastName/EmpId columns, and fill them with
#   the admin/admin login used before adding an employee
# - Active sheet/selection moves to the "Add Employee" sheet

$wb = $excel.ActiveWorkbook

# --- TestCases sheet: fill in the Result column for two more rows ---
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestCases.Range("G3").Value = "Fail"
$wsTestCases.Range("G7").Value = "Pass"

# --- Add Employee sheet: insert two new columns (UserName, Password) ---
$wsAddEmployee = $wb.Worksheets.Item("Add Employee")
$wsAddEmployee.Range("D1:E1").EntireColumn.Insert() | Out-Null

$wsAddEmployee.Range("D1").Value = "UserName"
$wsAddEmployee.Range("E1").Value = "Password"

$wsAddEmployee.Range("D2").Value = "admin"
$wsAddEmployee.Range("E2").Value = "admin"
$wsAddEmployee.Range("D3").Value = "admin"
$wsAddEmployee.Range("E3").Value = "admin"
$wsAddEmployee.Range("D4").Value = "admin"
$wsAddEmployee.Range("E4").Value = "admin"

$wsAddEmployee.Columns.Item(4).ColumnWidth = 11.592447916666666
$wsAddEmployee.Columns.Item(5).ColumnWidth = 11.592447916666666

# --- Selection / active sheet bookkeeping ---
$wsTestCases.Range("H4").Select() | Out-Null
$wsAddEmployee.Activate() | Out-Null
$wsAddEmployee.Range("E18").Select() | Out-Null
